# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sun Sep 17 20:47:38 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells B2:E51 hold text values (coin name / link / formatted price / %),
# never numbers -- some of the new price strings (e.g. '0.499') look like
# plain numbers, so force Text format first or Excel would silently store
# them as numeric values instead of strings.
$ws.Range("B2:E51").NumberFormat = "@"

$updates = @(
    @('D2', '26.729.95'),
    @('E2', '  -0.14%  '),
    @('D3', '1.638.45'),
    @('E3', '  -0.60%  '),
    @('E4', '  +0.24%  '),
    @('D5', '217.76'),
    @('E5', '  +0.56%  '),
    @('D6', '0.499'),
    @('E6', '  -1.42%  '),
    @('E7', '  +0.27%  '),
    @('D8', '0.248'),
    @('E8', '  -1.37%  '),
    @('D9', '0.0619'),
    @('E9', '  -1.34%  '),
    @('D10', '18.93'),
    @('E10', '  -1.64%  '),
    @('D11', '0.0844'),
    @('E11', '  +0.06%  '),
    @('D12', '1.871.07'),
    @('E12', '  -0.39%  '),
    @('D13', '1.641.34'),
    @('E13', '  -0.28%  '),
    @('D14', '4.11'),
    @('E14', '  -2.07%  '),
    @('D15', '0.521'),
    @('E15', '  -2.24%  '),
    @('D16', '64.01'),
    @('E16', '  -2.07%  '),
    @('D17', '26.751.23'),
    @('E17', '  -0.10%  '),
    @('D18', '0.0₃0720'),
    @('E18', '  -3.18%  '),
    @('E19', '  +0.13%  '),
    @('D20', '209.23'),
    @('E20', '  -3.89%  '),
    @('D21', '4.31'),
    @('E21', '  -1.27%  '),
    @('D22', '6.15'),
    @('E22', '  -1.77%  '),
    @('E23', '  -5.24%  '),
    @('D24', '9.17'),
    @('E24', '  -3.38%  '),
    @('D25', '146.85'),
    @('E25', '  -0.05%  '),
    @('E26', '  -0.01%  '),
    @('D27', '0.117'),
    @('E27', '  -2.52%  '),
    @('D28', '7.02'),
    @('E28', '  -2.21%  '),
    @('D29', '15.46'),
    @('E29', '  -1.91%  '),
    @('D30', '0.0499'),
    @('E30', '  -3.86%  '),
    @('D31', '1.18'),
    @('E31', '  +0.36%  '),
    @('D32', '3.32'),
    @('E32', '  -1.16%  '),
    @('D33', '2.94'),
    @('E33', '  -2.27%  '),
    @('D34', '1.266.36'),
    @('E34', '  -1.13%  '),
    @('B35', 'HuobiToken'),
    @('C35', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @('D35', '2.45'),
    @('E35', '  +0.09%  '),
    @('B36', 'LidoDAOToken'),
    @('C36', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'),
    @('D36', '1.51'),
    @('E36', '  -2.01%  '),
    @('D37', '0.0173'),
    @('E37', '  -3.44%  '),
    @('D38', '0.521'),
    @('E38', '  -3.17%  '),
    @('B39', 'PaxDollar'),
    @('C39', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'),
    @('D39', '1.01'),
    @('E39', '  +0.14%  '),
    @('B40', 'ARBITRUM'),
    @('C40', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'),
    @('D40', '0.798'),
    @('E40', '  -3.86%  '),
    @('D41', '0.799'),
    @('E41', '  -2.10%  '),
    @('D42', '2.18'),
    @('E42', '  -2.90%  '),
    @('D43', '1.781.62'),
    @('E43', '  -0.42%  '),
    @('D44', '5.25'),
    @('E44', '  -3.54%  '),
    @('D45', '91.00'),
    @('E45', '  -1.11%  '),
    @('D46', '59.83'),
    @('E46', '  -0.01%  '),
    @('D47', '1.57'),
    @('E47', '  -2.47%  '),
    @('E48', '  +0.81%  '),
    @('B49', 'USDD'),
    @('C49', 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'),
    @('D49', '1.01'),
    @('E49', '  +0.38%  '),
    @('B50', 'Mantle'),
    @('C50', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D50', '0.407'),
    @('E50', '  -0.26%  '),
    @('B51', 'EnergySwap'),
    @('C51', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @('D51', '7.49'),
    @('E51', '  -3.53%  ')
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Drop the temporary Text number-format override so the cells fall back
# to the workbook's normal (unstyled) look, matching the original sheet.
$ws.Range("B2:E51").Style = "Normal"
